# Insert two new data rows into the "Limón" price sheet at row 282,
# pushing the existing rows 282-349 down to 284-351, then populate the
# two new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 282 (shifts 282:349 -> 284:351)
$ws.Rows("282:283").Insert()

# Boilerplate values shared by every data row in this sheet
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"

# --- New row 282 ---
$ws.Range("A282").Value2 = $mercadoId
$ws.Range("B282").Value2 = $mercado
$ws.Range("C282").Value2 = $region
$ws.Range("D282").Value2 = 44932
$ws.Range("E282").Value2 = $codreg
$ws.Range("F282").Value2 = $tipo
$ws.Range("G282").Value2 = $productoId
$ws.Range("H282").Value2 = $producto
$ws.Range("I282").Value2 = $categoriaId
$ws.Range("J282").Value2 = $categoria
$ws.Range("K282").Value2 = "Sutil De Gase"
$ws.Range("L282").Value2 = "Primera"
$ws.Range("M282").Value2 = 250
$ws.Range("N282").Value2 = 23000
$ws.Range("O282").Value2 = 24000
$ws.Range("P282").Value2 = 23500
$ws.Range("Q282").Value2 = "`$/caja 18 kilos"
$ws.Range("R282").Value2 = "Perú"
$ws.Range("S282").Value2 = 1306
$ws.Range("T282").Value2 = 18

# --- New row 283 ---
$ws.Range("A283").Value2 = $mercadoId
$ws.Range("B283").Value2 = $mercado
$ws.Range("C283").Value2 = $region
$ws.Range("D283").Value2 = 44932
$ws.Range("E283").Value2 = $codreg
$ws.Range("F283").Value2 = $tipo
$ws.Range("G283").Value2 = $productoId
$ws.Range("H283").Value2 = $producto
$ws.Range("I283").Value2 = $categoriaId
$ws.Range("J283").Value2 = $categoria
$ws.Range("K283").Value2 = "Tahití"
$ws.Range("L283").Value2 = "Primera"
$ws.Range("M283").Value2 = 250
$ws.Range("N283").Value2 = 30000
$ws.Range("O283").Value2 = 31000
$ws.Range("P283").Value2 = 30500
$ws.Range("Q283").Value2 = "`$/caja 24 kilos"
$ws.Range("R283").Value2 = "Perú"
$ws.Range("S283").Value2 = 1271
$ws.Range("T283").Value2 = 24
